{"js": "// Apply the documented text edits to the BAFDR software README.\n// Each change is performed as an exact-text search within the document\n// body followed by a \"Replace\" insertText, which keeps the surrounding\n// run/formatting structure intact while updating the wording.\n\nasync function replaceOnce(body, oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) \"...execute each module... website.  Then, download the files to a\n//     directory/folder of your choosing.\" -> mention github.\nawait replaceOnce(\n  body,\n  \" Then, download the files to a directory/folder of your choosing.\",\n  \" Then, download the files from github to a directory/folder of your choosing.\"\n);\n\n// 2) \"Then you can execute the notebook.\" -> add lead-in sentence.\nawait replaceOnce(\n  body,\n  \"Then you can execute the\",\n  \"Update this as needed, then you can execute the\"\n);\n\n// 3) \"You can use excel or similar to edit the BAFDR...\" -> add \"program\".\nawait replaceOnce(\n  body,\n  \"You can use excel or similar to edit the BAFDR\",\n  \"You can use excel or similar program to edit the BAFDR\"\n);\n\n// 4) Replace the \"refer to the BAFDR_S_input.txt...\" guidance with the new\n//    sentences about headers and unlimited case counts.\nawait replaceOnce(\n  body,\n  \"The input order is the same as the single case version but transposed into a row, since there is not currently text guidance you must refer to the BAFDR_S_input.txt in case of any question on the order of parameters.  \",\n  \"The input order is the same as the single case version but transposed into a row.  The first row is the headers with variable names. There is no limit set on the number of cases.  \"\n);\n\n// 5) Parenthesize \"on line 14 in this case\".\nawait replaceOnce(\n  body,\n  \" The ipynb file again needs to be edited to correspond to this directory on line 14 in this case. \",\n  \" The ipynb file again needs to be edited to correspond to this directory (on line 14 in this case). \"\n);\n\n// 6) \"with no plots is very quick\" -> \"without plots is quick\".\nawait replaceOnce(\n  body,\n  \" No plots are output for the multi-case runs, just the csv file.  The execution with no plots is very quick.  \",\n  \" No plots are output for the multi-case runs, just the csv file.  The execution without plots is quick.  \"\n);\n", "ps1": "# Apply the documented text edits to the BAFDR software README using\n# Word's Find/Replace (COM object model). Each call targets one exact\n# phrase so only the intended text changes.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\n# 1) \"...website.  Then, download the files to a directory/folder of your\n#     choosing.\" -> mention github.\nReplace-Text \" Then, download the files to a directory/folder of your choosing.\" \" Then, download the files from github to a directory/folder of your choosing.\"\n\n# 2) \"Then you can execute the notebook.\" -> add lead-in sentence.\nReplace-Text \"Then you can execute the\" \"Update this as needed, then you can execute the\"\n\n# 3) \"You can use excel or similar to edit the BAFDR...\" -> add \"program\".\nReplace-Text \"You can use excel or similar to edit the BAFDR\" \"You can use excel or similar program to edit the BAFDR\"\n\n# 4) Replace the \"refer to the BAFDR_S_input.txt...\" guidance with the new\n#    sentences about headers and unlimited case counts.\nReplace-Text \"The input order is the same as the single case version but transposed into a row, since there is not currently text guidance you must refer to the BAFDR_S_input.txt in case of any question on the order of parameters.  \" \"The input order is the same as the single case version but transposed into a row.  The first row is the headers with variable names. There is no limit set on the number of cases.  \"\n\n# 5) Parenthesize \"on line 14 in this case\".\nReplace-Text \" The ipynb file again needs to be edited to correspond to this directory on line 14 in this case. \" \" The ipynb file again needs to be edited to correspond to this directory (on line 14 in this case). \"\n\n# 6) \"with no plots is very quick\" -> \"without plots is quick\".\nReplace-Text \" No plots are output for the multi-case runs, just the csv file.  The execution with no plots is very quick.  \" \" No plots are output for the multi-case runs, just the csv file.  The execution without plots is quick.  \"\n"}
